$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 56, pushing existing rows 56..176 down to 57..177.
$ws.Rows.Item(56).Insert()

# Fill in the new weekly record on the newly inserted row 56.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across the whole table.
$ws.Range("A56").Value = 4
$ws.Range("B56").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C56").Value = "Los Lagos"
$ws.Range("E56").Value = 10
$ws.Range("F56").Value = 100112003
$ws.Range("G56").Value = "Ajo"
$ws.Range("H56").Value = "Chino"
$ws.Range("I56").Value = "Primera"
$ws.Range("N56").Value = "$/caja 10 kilos"
$ws.Range("O56").Value = "China"
$ws.Range("Q56").Value = 10
$ws.Range("R56").Value = "Hortaliza"

# New record-specific values.
$ws.Range("D56").Value = 44519
$ws.Range("J56").Value = 200
$ws.Range("K56").Value = 22000
$ws.Range("L56").Value = 23000
$ws.Range("M56").Value = 22500
$ws.Range("P56").Value = 2250
